# Applies the "Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig"
# content refresh to the StructureDefinition workbook:
#   - rebrand IBM/Alvearie -> LinuxForHealth in the canonical URL + Publisher
#   - bump Version 7.0.0 -> 8.0.0 and the generation Date
#   - clear the stale ele-1/ext-1 Constraint(s) note on the root "Extension" row
#     (the IG publisher no longer emits it there in the regenerated sheet)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-allowed-amount-total"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-allowed-amount-total"
$wsElements.Range("AI2").Value = ""
